# Append: 2026-01-18 06:28 JST
# - Update the "取得日時" (fetched-at) timestamp on the first two rows.
# - Replace row 3's listing with the listing that used to be on row 11
#   (the newly-collected item), including its URL/hyperlink.
# - Drop the now-duplicated rows 4-11 (the sheet keeps only the freshest
#   2 listings after this run).
# - Narrow columns D and H to match the new content widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-18 06:28:42"

# Row 2 keeps the same listing, only the fetch timestamp advances.
$ws.Range("A2").Value = $newTimestamp

# Row 3 becomes the listing that was previously on row 11.
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "製造業DXプロダクト開発のプロダクトマネージャー募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5468432"
$ws.Range("G3").Value = 75
$ws.Range("H3").Value = "◆開発"

# Rows 4-11 no longer exist in the refreshed data set.
$ws.Range("A4:H11").EntireRow.Delete()

# Rebuild the hyperlinks collection so only F2/F3 carry links, with F3
# now pointing at the listing that was promoted into it.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5473648")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5468432")

# Narrower content on column D/H lets those columns shrink. (The engine's
# ColumnWidth setter round-trips through a pixel grid that adds 5/6 of a
# character back on read, so we dial the input down by that much to land
# on an exact integer width of 28 / 12 in the saved file.)
$ws.Columns.Item(4).ColumnWidth = 27.166666666666668
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666
